# Add a new trailing data row (row 73) to each of the four worksheets,
# mirroring the existing row layout (time, length/id/checksum hex blobs,
# and their decoded numeric counterparts).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        SheetName = "ROW35-FE-LIFTER"
        A = 45760.89821945602
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x66"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 13
    },
    @{
        SheetName = "ROW35-MID-LIFTER"
        A = 45760.7550420949
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x66"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 14
    },
    @{
        SheetName = "ROW02-FE-LIFTER"
        A = 45760.89562313657
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x66"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 358
        I = 3
    },
    @{
        SheetName = "ROW02-MID-LIFTER"
        A = 45760.95647961806
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x66"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 358
        I = 3
    }
)

foreach ($rowData in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($rowData.SheetName)
    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $rowData.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F
    $ws.Cells.Item($newRow, 7).Value = $rowData.G
    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}
